$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3111703333333333
$ws.Range("H2").Value = 0.933511
$ws.Range("I2").Value = 0.1268479703884478
$ws.Range("J2").Value = 0.1268479703884478
$ws.Range("M2").Value = 7.407905
$ws.Range("N2").Value = 22.223715
$ws.Range("O2").Value = 0.1577242380174723
$ws.Range("P2").Value = 0.1577242380174723
$ws.Range("Q2").Value = 2.305120268151666
$ws.Range("R2").Value = 20.746082413365
$ws.Range("S2").Value = 0.02000699947358081
$ws.Range("T2").Value = 0.02000699947358082

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3111703333333333
$ws.Range("H3").Value = 0.933511
$ws.Range("I3").Value = 0.1268479703884478
$ws.Range("J3").Value = 0.1268479703884478
$ws.Range("O3").Value = 0.3510414535684271
$ws.Range("P3").Value = 0.3510414535684271
$ws.Range("Q3").Value = 5.130427509133778
$ws.Range("R3").Value = 46.173847582204
$ws.Range("S3").Value = 0.04452889590736551
$ws.Range("T3").Value = 0.04452889590736551

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3111703333333333
$ws.Range("H4").Value = 0.933511
$ws.Range("I4").Value = 0.1268479703884478
$ws.Range("J4").Value = 0.1268479703884478
$ws.Range("M4").Value = 5.464566666666666
$ws.Range("N4").Value = 16.3937
$ws.Range("O4").Value = 0.1163479571613943
$ws.Range("P4").Value = 0.1163479571613943
$ws.Range("Q4").Value = 1.700411031188889
$ws.Range("R4").Value = 15.3036992807
$ws.Range("S4").Value = 0.01475850222476493
$ws.Range("T4").Value = 0.01475850222476493

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3111703333333333
$ws.Range("H5").Value = 0.933511
$ws.Range("I5").Value = 0.1268479703884478
$ws.Range("J5").Value = 0.1268479703884478
$ws.Range("M5").Value = 17.60745533333333
$ws.Range("N5").Value = 52.822366
$ws.Range("O5").Value = 0.3748863512527063
$ws.Range("P5").Value = 0.3748863512527063
$ws.Range("Q5").Value = 5.478917745225111
$ws.Range("R5").Value = 49.310259707026
$ws.Range("S5").Value = 0.04755357278273652
$ws.Range("T5").Value = 0.04755357278273652

# Row 6
$ws.Range("G6").Value = 2.141926333333334
$ws.Range("H6").Value = 6.425779
$ws.Range("I6").Value = 0.8731520296115521
$ws.Range("J6").Value = 0.8731520296115521
$ws.Range("M6").Value = 7.407905
$ws.Range("N6").Value = 22.223715
$ws.Range("O6").Value = 0.1577242380174723
$ws.Range("P6").Value = 0.1577242380174723
$ws.Range("Q6").Value = 15.86718679433167
$ws.Range("R6").Value = 142.804681148985
$ws.Range("S6").Value = 0.1377172385438915
$ws.Range("T6").Value = 0.1377172385438915

# Row 7
$ws.Range("G7").Value = 2.141926333333334
$ws.Range("H7").Value = 6.425779
$ws.Range("I7").Value = 0.8731520296115521
$ws.Range("J7").Value = 0.8731520296115521
$ws.Range("O7").Value = 0.3510414535684271
$ws.Range("P7").Value = 0.3510414535684271
$ws.Range("Q7").Value = 35.31505611526178
$ws.Range("R7").Value = 317.835505037356
$ws.Range("S7").Value = 0.3065125576610616
$ws.Range("T7").Value = 0.3065125576610616

# Row 8
$ws.Range("G8").Value = 2.141926333333334
$ws.Range("H8").Value = 6.425779
$ws.Range("I8").Value = 0.8731520296115521
$ws.Range("J8").Value = 0.8731520296115521
$ws.Range("M8").Value = 5.464566666666666
$ws.Range("N8").Value = 16.3937
$ws.Range("O8").Value = 0.1163479571613943
$ws.Range("P8").Value = 0.1163479571613943
$ws.Range("Q8").Value = 11.70469924358889
$ws.Range("R8").Value = 105.3422931923
$ws.Range("S8").Value = 0.1015894549366293
$ws.Range("T8").Value = 0.1015894549366293

# Row 9
$ws.Range("G9").Value = 2.141926333333334
$ws.Range("H9").Value = 6.425779
$ws.Range("I9").Value = 0.8731520296115521
$ws.Range("J9").Value = 0.8731520296115521
$ws.Range("M9").Value = 17.60745533333333
$ws.Range("N9").Value = 52.822366
$ws.Range("O9").Value = 0.3748863512527063
$ws.Range("P9").Value = 0.3748863512527063
$ws.Range("Q9").Value = 37.71387224145712
$ws.Range("R9").Value = 339.424850173114
$ws.Range("S9").Value = 0.3273327784699697
$ws.Range("T9").Value = 0.3273327784699698

